$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 700
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 3000
$ws.Range("N29").Value = -3562
$ws.Range("H38").Value = 2688276.2
$ws.Range("I38").Value = 2688276.2
$ws.Range("K38").Value = 8064828.600000001
$ws.Range("M38").Value = -8064456.600000001
$ws.Range("H43").Value = 2740.2
$ws.Range("I43").Value = 5500.5
$ws.Range("J43").Value = 900
$ws.Range("K43").Value = 5500.5
$ws.Range("L43").Value = 900
$ws.Range("M43").Value = -5431.5
$ws.Range("N43").Value = -1038
$ws.Range("H58").Value = 817930.3
$ws.Range("I58").Value = 1089579.2
$ws.Range("J58").Value = 2983.3333
$ws.Range("K58").Value = 3268737.6
$ws.Range("L58").Value = 8949.999899999999
$ws.Range("M58").Value = -3268587.6
$ws.Range("N58").Value = -9249.999899999999
$ws.Range("H63").Value = 35000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 35000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 35000
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -36248
$ws.Range("H66").Value = 35000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 35000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 105000
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -111240
$ws.Range("H86").Value = 51701.6
$ws.Range("I86").Value = 51701.6
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 51701.6
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -50578.6
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 51701.6
$ws.Range("I89").Value = 51701.6
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 258508
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -252892
$ws.Range("N89").ClearContents()
$ws.Range("H125").Value = 1831.3462
$ws.Range("I125").Value = 1984.7778
$ws.Range("J125").Value = 1750.1177
$ws.Range("K125").Value = 17863.0002
$ws.Range("L125").Value = 15751.0593
$ws.Range("M125").Value = -15403.0002
$ws.Range("N125").Value = -20671.0593
$ws.Range("H132").Value = 16680321
$ws.Range("I132").Value = 17871630
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 53614890
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -53612360
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26575.465
$ws.Range("I32").Value = 5899.9243
$ws.Range("J32").Value = 221516.28
$ws.Range("K32").Value = 5899.9243
$ws.Range("L32").Value = 221516.28
$ws.Range("M32").Value = -5612.9243
$ws.Range("N32").Value = -222090.28
$ws.Range("H45").Value = 77752.92
$ws.Range("I45").Value = 200831.2
$ws.Range("J45").Value = 829
$ws.Range("K45").Value = 200831.2
$ws.Range("L45").Value = 829
$ws.Range("M45").Value = -200454.2
$ws.Range("N45").Value = -1583
$ws.Range("H61").Value = 1855.238
$ws.Range("I61").Value = 1740
$ws.Range("J61").Value = 2950
$ws.Range("K61").Value = 1740
$ws.Range("L61").Value = 2950
$ws.Range("M61").Value = -1528
$ws.Range("N61").Value = -3374
$ws.Range("H122").Value = 1293.091
$ws.Range("I122").Value = 1255.68
$ws.Range("J122").Value = 1410
$ws.Range("K122").Value = 3767.04
$ws.Range("L122").Value = 4230
$ws.Range("M122").Value = -1317.04
$ws.Range("N122").Value = -9130
$ws.Range("H136").Value = 1855.238
$ws.Range("I136").Value = 1740
$ws.Range("J136").Value = 2950
$ws.Range("K136").Value = 5220
$ws.Range("L136").Value = 8850
$ws.Range("M136").Value = -2670
$ws.Range("N136").Value = -13950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2204.875
$ws.Range("I99").Value = 1995.6
$ws.Range("J99").Value = 2300
$ws.Range("K99").Value = 1995.6
$ws.Range("L99").Value = 2300
$ws.Range("M99").Value = -497.5999999999999
$ws.Range("N99").Value = -5296
$ws.Range("H134").Value = 2031.1471
$ws.Range("I134").Value = 2157.2222
$ws.Range("J134").Value = 1544.8572
$ws.Range("K134").Value = 6471.6666
$ws.Range("L134").Value = 4634.571599999999
$ws.Range("M134").Value = -3936.6666
$ws.Range("N134").Value = -9704.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 11332.667
$ws.Range("I45").Value = 9499
$ws.Range("K45").Value = 9499
$ws.Range("M45").Value = -8906
$ws.Range("H58").Value = 1872.0769
$ws.Range("I58").Value = 1932.5
$ws.Range("J58").Value = 1820.2858
$ws.Range("K58").Value = 1932.5
$ws.Range("L58").Value = 1820.2858
$ws.Range("M58").Value = -1729.5
$ws.Range("N58").Value = -2226.2858
$ws.Range("H99").Value = 8910.883
$ws.Range("I99").Value = 3810.1428
$ws.Range("J99").Value = 12481.4
$ws.Range("K99").Value = 3810.1428
$ws.Range("L99").Value = 12481.4
$ws.Range("M99").Value = -2312.1428
$ws.Range("N99").Value = -15477.4
$ws.Range("H126").Value = 8910.883
$ws.Range("I126").Value = 3810.1428
$ws.Range("J126").Value = 12481.4
$ws.Range("K126").Value = 11430.4284
$ws.Range("L126").Value = 37444.2
$ws.Range("M126").Value = -8960.428400000001
$ws.Range("N126").Value = -42384.2
$ws.Range("H134").Value = 1264
$ws.Range("I134").Value = 863.2
$ws.Range("J134").Value = 2600
$ws.Range("K134").Value = 2589.6
$ws.Range("L134").Value = 7800
$ws.Range("M134").Value = -54.60000000000036
$ws.Range("N134").Value = -12870
$ws.Range("H136").Value = 1872.0769
$ws.Range("I136").Value = 1932.5
$ws.Range("J136").Value = 1820.2858
$ws.Range("K136").Value = 5797.5
$ws.Range("L136").Value = 5460.857400000001
$ws.Range("M136").Value = -3247.5
$ws.Range("N136").Value = -10560.8574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 2449.75
$ws.Range("J58").Value = 2966.3333
$ws.Range("L58").Value = 8898.999899999999
$ws.Range("N58").Value = -9154.999899999999
$ws.Range("H59").Value = 890.8
$ws.Range("I59").Value = 818
$ws.Range("J59").Value = 1000
$ws.Range("K59").Value = 2454
$ws.Range("L59").Value = 3000
$ws.Range("M59").Value = -1914
$ws.Range("N59").Value = -4080
$ws.Range("H131").Value = 10902.76
$ws.Range("I131").Value = 1211.6666
$ws.Range("J131").Value = 11521.341
$ws.Range("K131").Value = 3634.9998
$ws.Range("L131").Value = 34564.023
$ws.Range("M131").Value = 1405.0002
$ws.Range("N131").Value = -44644.023

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 5
$ws.Range("I18").Value = 5
$ws.Range("K18").Value = 5
$ws.Range("M18").Value = 288
$ws.Range("H80").Value = 100105350
$ws.Range("J80").Value = 3312
$ws.Range("L80").Value = 3312
$ws.Range("N80").Value = -5308
$ws.Range("H83").Value = 100105350
$ws.Range("J83").Value = 3312
$ws.Range("L83").Value = 16560
$ws.Range("N83").Value = -26544
$ws.Range("H102").Value = 2971.4
$ws.Range("I102").Value = 2032.3334
$ws.Range("J102").Value = 4380
$ws.Range("K102").Value = 2032.3334
$ws.Range("L102").Value = 4380
$ws.Range("M102").Value = -410.3334
$ws.Range("N102").Value = -7624
$ws.Range("H122").Value = 987.913
$ws.Range("I122").Value = 910.1667
$ws.Range("J122").Value = 1072.7273
$ws.Range("K122").Value = 2730.5001
$ws.Range("L122").Value = 3218.1819
$ws.Range("M122").Value = -280.5001000000002
$ws.Range("N122").Value = -8118.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1368
$ws.Range("I122").Value = 1221
$ws.Range("J122").Value = 1662
$ws.Range("K122").Value = 3663
$ws.Range("L122").Value = 4986
$ws.Range("M122").Value = -1213
$ws.Range("N122").Value = -9886
